$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "oil_economics" worksheet, positioned right after
#    "pel_economics" (so the sheet order becomes: gen_economics,
#    gas_economics, el_economics, pel_economics, oil_economics,
#    dev_economics, comp_economics, ep_table, further_parameters).
# ---------------------------------------------------------------------------
$pel = $wb.Worksheets.Item("pel_economics")
$oil = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $pel)
$oil.Name = "oil_economics"

# Header row, mirrors the layout used on the other *_economics sheets
# (Name / fixed fees / variable price / CO2 emissions / Source).
$oil.Cells.Item(1,1).Value = "Name"
$oil.Cells.Item(1,2).Value = "fixed fees [Euro/a]"
$oil.Cells.Item(1,3).Value = "variable price [Euro/kWh]"
$oil.Cells.Item(1,4).Value = "CO2 emissions [kgCO2/kWh]"
$oil.Cells.Item(1,5).Value = "Source"

# Data row for the (new) old oil boiler tariff "oil_sta".
$oil.Cells.Item(2,1).Value = "oil_sta"
$oil.Cells.Item(2,2).Value = "[0 - 9999: 1]"
$oil.Cells.Item(2,3).Value = "[0 - 9999: 0.048]"
$oil.Cells.Item(2,4).Value = 0.025
$oil.Range("E2").Style = "Link"

$oil.Range("B38").Select()

# ---------------------------------------------------------------------------
# 2) Insert a new row into "gen_economics" for the oil price-change factor,
#    right above the EEX compensation row.
# ---------------------------------------------------------------------------
$gen = $wb.Worksheets.Item("gen_economics")
$gen.Rows.Item(8).Insert()
$gen.Cells.Item(8,1).Value = "prChange_oil"
$gen.Cells.Item(8,2).Value = 1.001
$gen.Cells.Item(8,3).Value = "-"
$gen.Cells.Item(8,4).Value = "Price change factors per year for oil"

# ---------------------------------------------------------------------------
# 3) Restore a sensible view state: pel_economics keeps a block selection,
#    gen_economics becomes the active sheet/cell.
# ---------------------------------------------------------------------------
$pel.Activate()
$pel.Range("A1:E2").Select()

$gen.Activate()
$gen.Range("B17").Select()
